$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.415.12'
$ws.Range("E2").Value = '  -2.91%  '

$ws.Range("D3").Value = '1.773.40'
$ws.Range("E3").Value = '  -1.98%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.01'
$ws.Range("E6").Value = '  -1.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4264'
$ws.Range("E7").Value = '  +1.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3617'
$ws.Range("E8").Value = '  +1.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07131'
$ws.Range("E9").Value = '  +0.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8385'
$ws.Range("E10").Value = '  -1.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.39'
$ws.Range("E11").Value = '  +0.88%  '

$ws.Range("D12").Value = '1.802.44'
$ws.Range("E12").Value = '  +0.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.441'
$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.248'
$ws.Range("E14").Value = '  -1.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06908'
$ws.Range("E15").Value = '  +0.70%  '

$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.96'
$ws.Range("E17").Value = '  -2.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008692'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.93'
$ws.Range("E20").Value = '  -1.04%  '

$ws.Range("D21").Value = '26.437.17'
$ws.Range("E21").Value = '  -2.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.097'
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.06'

$ws.Range("D24").Value = '2.018.44'
$ws.Range("E24").Value = '  +0.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.28'
$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.808'
$ws.Range("E26").Value = '  -8.61%  '

$ws.Range("E27").Value = '  -0.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.065'
$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.86'
$ws.Range("E29").Value = '  +0.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.776'
$ws.Range("E30").Value = '  +4.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08880'
$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7254'
$ws.Range("E32").Value = '  -2.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.114'
$ws.Range("E33").Value = '  +1.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.318'
$ws.Range("E34").Value = '  -2.83%  '

$ws.Range("E35").Value = '  +0.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.748'
$ws.Range("E36").Value = '  -5.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.101'
$ws.Range("E37").Value = '  +3.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05138'
$ws.Range("E38").Value = '  -0.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01882'
$ws.Range("E39").Value = '  -0.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1613'
$ws.Range("E40").Value = '  -1.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4905'
$ws.Range("E41").Value = '  -1.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.593'
$ws.Range("E42").Value = '  -4.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.347'
$ws.Range("E43").Value = '  +1.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.962'
$ws.Range("E44").Value = '  -2.26%  '

$ws.Range("E45").Value = '  -0.40%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.15'
$ws.Range("E47").Value = '  -0.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.630'
$ws.Range("E48").Value = '  +2.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06180'
$ws.Range("E49").Value = '  -3.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4431'
$ws.Range("E50").Value = '  -2.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.717'
$ws.Range("E51").Value = '  +0.69%  '
